$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 295.83334
$ws.Range("I33").Value = 227.3077
$ws.Range("K33").Value = 227.3077
$ws.Range("M33").Value = 1.692299999999989
$ws.Range("H80").Value = 1537.6666
$ws.Range("I80").Value = 469.9091
$ws.Range("J80").Value = 2271.75
$ws.Range("K80").Value = 1409.7273
$ws.Range("L80").Value = 6815.25
$ws.Range("M80").Value = -411.7273
$ws.Range("N80").Value = -8811.25
$ws.Range("H83").Value = 1537.6666
$ws.Range("I83").Value = 469.9091
$ws.Range("J83").Value = 2271.75
$ws.Range("K83").Value = 4229.1819
$ws.Range("L83").Value = 20445.75
$ws.Range("M83").Value = 762.8180999999995
$ws.Range("N83").Value = -30429.75
$ws.Range("H94").Value = 1561.5555
$ws.Range("I94").Value = 1474
$ws.Range("J94").Value = 1868
$ws.Range("K94").Value = 1474
$ws.Range("L94").Value = 1868
$ws.Range("M94").Value = -1023
$ws.Range("N94").Value = -2770
$ws.Range("H100").Value = 6320.92
$ws.Range("I100").Value = 5497.727
$ws.Range("J100").Value = 6967.7144
$ws.Range("K100").Value = 5497.727
$ws.Range("L100").Value = 6967.7144
$ws.Range("M100").Value = -4956.727
$ws.Range("N100").Value = -8049.7144
$ws.Range("H127").Value = 6054.5557
$ws.Range("I127").Value = 1686.375
$ws.Range("K127").Value = 5059.125
$ws.Range("M127").Value = -99.125
$ws.Range("H137").Value = 5715.154
$ws.Range("I137").Value = 5662.1875
$ws.Range("K137").Value = 16986.5625
$ws.Range("M137").Value = -14436.5625
$ws.Range("H138").Value = 3042.8604
$ws.Range("J138").Value = 3700.1667
$ws.Range("L138").Value = 11100.5001
$ws.Range("N138").Value = -21380.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8234.177
$ws.Range("I61").Value = 1485.091
$ws.Range("K61").Value = 1485.091
$ws.Range("M61").Value = -1273.091
$ws.Range("H63").Value = 7456.609
$ws.Range("I63").Value = 5958.6665
$ws.Range("K63").Value = 5958.6665
$ws.Range("M63").Value = -5272.6665
$ws.Range("H66").Value = 7456.609
$ws.Range("I66").Value = 5958.6665
$ws.Range("K66").Value = 29793.3325
$ws.Range("M66").Value = -26361.3325
$ws.Range("H74").Value = 3990.2
$ws.Range("I74").Value = 3632.3572
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 3632.3572
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -2758.3572
$ws.Range("N74").Value = -10748
$ws.Range("H77").Value = 3990.2
$ws.Range("I77").Value = 3632.3572
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 18161.786
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -13793.786
$ws.Range("N77").Value = -53736
$ws.Range("H102").Value = 33343756
$ws.Range("I102").Value = 6777.5
$ws.Range("K102").Value = 6777.5
$ws.Range("M102").Value = -5155.5
$ws.Range("H136").Value = 8234.177
$ws.Range("I136").Value = 1485.091
$ws.Range("K136").Value = 4455.272999999999
$ws.Range("M136").Value = -1905.272999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5968.6875
$ws.Range("I20").Value = 4472.222
$ws.Range("J20").Value = 7892.7144
$ws.Range("K20").Value = 4472.222
$ws.Range("L20").Value = 7892.7144
$ws.Range("M20").Value = -4225.222
$ws.Range("N20").Value = -8386.714400000001
$ws.Range("H86").Value = 22774356
$ws.Range("I86").Value = 31313990
$ws.Range("J86").Value = 1999.1666
$ws.Range("K86").Value = 31313990
$ws.Range("L86").Value = 1999.1666
$ws.Range("M86").Value = -31312867
$ws.Range("N86").Value = -4245.1666
$ws.Range("H89").Value = 22774356
$ws.Range("I89").Value = 31313990
$ws.Range("J89").Value = 1999.1666
$ws.Range("K89").Value = 156569950
$ws.Range("L89").Value = 9995.833000000001
$ws.Range("M89").Value = -156564334
$ws.Range("N89").Value = -21227.833
$ws.Range("H99").Value = 32542.54
$ws.Range("I99").Value = 41436
$ws.Range("K99").Value = 41436
$ws.Range("M99").Value = -39938
$ws.Range("H105").Value = 2917.4
$ws.Range("I105").Value = 2305.4443
$ws.Range("J105").Value = 3418.0908
$ws.Range("K105").Value = 2305.4443
$ws.Range("L105").Value = 3418.0908
$ws.Range("M105").Value = -558.4443000000001
$ws.Range("N105").Value = -6912.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3890.0715
$ws.Range("I31").Value = 3120.25
$ws.Range("K31").Value = 3120.25
$ws.Range("M31").Value = -2825.25
$ws.Range("H34").Value = 3890.0715
$ws.Range("I34").Value = 3120.25
$ws.Range("K34").Value = 3120.25
$ws.Range("M34").Value = -2918.25
$ws.Range("H35").Value = 1121.3334
$ws.Range("I35").Value = 1121.3334
$ws.Range("K35").Value = 1121.3334
$ws.Range("M35").Value = -827.3334
$ws.Range("H97").Value = 43166.332
$ws.Range("J97").Value = 43166.332
$ws.Range("L97").Value = 43166.332
$ws.Range("N97").Value = -45148.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 9748
$ws.Range("J32").Value = 9748
$ws.Range("L32").Value = 29244
$ws.Range("N32").Value = -29810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 2000
$ws.Range("I18").Value = 2000
$ws.Range("K18").Value = 2000
$ws.Range("M18").Value = -1828
$ws.Range("H22").Value = 18183016
$ws.Range("I22").Value = 30303696
$ws.Range("K22").Value = 30303696
$ws.Range("M22").Value = -30303401
$ws.Range("H27").Value = 18183016
$ws.Range("I27").Value = 30303696
$ws.Range("K27").Value = 30303696
$ws.Range("M27").Value = -30303589
$ws.Range("H32").Value = 1439.7142
$ws.Range("I32").Value = 1439.7142
$ws.Range("K32").Value = 1439.7142
$ws.Range("M32").Value = -1122.7142
$ws.Range("H61").Value = 2129.6667
$ws.Range("I61").Value = 2096
$ws.Range("K61").Value = 2096
$ws.Range("M61").Value = -1894
$ws.Range("H100").Value = 335333
$ws.Range("I100").Value = 335333
$ws.Range("K100").Value = 335333
$ws.Range("M100").Value = -334792
$ws.Range("H113").Value = 2129.6667
$ws.Range("I113").Value = 2096
$ws.Range("K113").Value = 2096
$ws.Range("M113").Value = 74

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 13845.5
$ws.Range("J40").Value = 13845.5
$ws.Range("L40").Value = 13845.5
$ws.Range("N40").Value = -14143.5
$ws.Range("H96").Value = 8000
$ws.Range("I96").Value = 6000
$ws.Range("J96").Value = 10000
$ws.Range("K96").Value = 6000
$ws.Range("L96").Value = 10000
$ws.Range("M96").Value = -4627
$ws.Range("N96").Value = -12746
$ws.Range("H109").Value = 18600
$ws.Range("J109").Value = 18600
$ws.Range("L109").Value = 18600
$ws.Range("N109").Value = -21374
$ws.Range("H132").Value = 4182.9346
$ws.Range("I132").Value = 2976.4849
$ws.Range("J132").Value = 7245.4614
$ws.Range("K132").Value = 8929.4547
$ws.Range("L132").Value = 21736.3842
$ws.Range("M132").Value = -6399.4547
$ws.Range("N132").Value = -26796.3842
$ws.Range("H138").Value = 67500
$ws.Range("J138").Value = 67500
$ws.Range("L138").Value = 67500
$ws.Range("N138").Value = -77780
